# Example.xlsx update:
#  1. add new "date" / "dateArry" / "timestamparry" columns (T, U, V) to Sheet1
#     mirroring the existing date/tinydate/timestamp/utctime block (columns P-S):
#       row2 -> field name
#       row3 -> field type
#       row4 -> field description (wrapped text, same style as the rest of row 4)
#       row5 -> example value
#  2. widen the new columns and move the sheet view/selection onto them
#  3. rename the default cell style from "Normal" to the workbook's own
#     localized name ("常规")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- new columns: content, written in the same order the original
#     author entered it (keeps the shared-string table identical) --------
$ws.Range("T3").Value = "{k:date}"
$ws.Range("U2").Value = "dateArry"
$ws.Range("T4").Value = "date in object value"
$ws.Range("U4").Value = "date array"
$ws.Range("U5").Value = '["2018/01/01 23:59:59"]'
$ws.Range("U3").Value = "date[1]"
$ws.Range("V2").Value = "timestamparry"
$ws.Range("V3").Value = "timestamp[2]"
$ws.Range("V4").Value = "timestamp array"
$ws.Range("V5").Value = "[1529995094, 1529995116]"
$ws.Range("T5").Value = '{"k":"2018/01/01 00:59:59"}'
$ws.Range("T2").Value = "date"

$ws.Range("T4:V4").WrapText = $true

# --- column widths for the new columns ----------------------------------
$ws.Columns.Item(20).ColumnWidth = 26.000468749999987
$ws.Columns.Item(21).ColumnWidth = 16.500156249999957
$ws.Columns.Item(22).ColumnWidth = 17.500156249999957

# --- move the view / selection onto the newly added columns ------------
$ws.Range("W5").Select()
$excel.ActiveWindow.ScrollColumn = 11

# --- rename the default style to match the workbook locale -------------
$wb.Styles.Item(1).Name = "常规"

$wb.Save()
